$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 4842.4287
$ws.Cells.Item(62, 10).Value = 5333
$ws.Cells.Item(62, 12).Value = 5333
$ws.Cells.Item(62, 14).Value = -6581

$ws.Cells.Item(65, 8).Value = 4842.4287
$ws.Cells.Item(65, 10).Value = 5333
$ws.Cells.Item(65, 12).Value = 26665
$ws.Cells.Item(65, 14).Value = -32905

$ws.Cells.Item(113, 8).Value = 3387.4546
$ws.Cells.Item(113, 9).Value = 3497.75
$ws.Cells.Item(113, 10).Value = 3324.4285
$ws.Cells.Item(113, 11).Value = 3497.75
$ws.Cells.Item(113, 12).Value = 3324.4285
$ws.Cells.Item(113, 13).Value = -243.75
$ws.Cells.Item(113, 14).Value = -9832.4285

$ws.Cells.Item(118, 8).Value = 494.125
$ws.Cells.Item(118, 9).Value = 408.6
$ws.Cells.Item(118, 11).Value = 1225.8
$ws.Cells.Item(118, 13).Value = 431.1999999999998

$ws.Cells.Item(132, 8).Value = 3976.889
$ws.Cells.Item(132, 9).Value = 1422.5667
$ws.Cells.Item(132, 10).Value = 16748.5
$ws.Cells.Item(132, 11).Value = 4267.7001
$ws.Cells.Item(132, 12).Value = 50245.5
$ws.Cells.Item(132, 13).Value = -1737.7001
$ws.Cells.Item(132, 14).Value = -55305.5

$ws.Cells.Item(135, 8).Value = 7917.467
$ws.Cells.Item(135, 9).Value = 1692.125
$ws.Cells.Item(135, 11).Value = 15229.125
$ws.Cells.Item(135, 13).Value = -12694.125

$ws.Cells.Item(138, 8).Value = 5845.5317
$ws.Cells.Item(138, 9).Value = 1719.3889
$ws.Cells.Item(138, 10).Value = 7063.082
$ws.Cells.Item(138, 11).Value = 5158.1667
$ws.Cells.Item(138, 12).Value = 21189.246
$ws.Cells.Item(138, 13).Value = -18.16669999999976
$ws.Cells.Item(138, 14).Value = -31469.246

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 693.2593000000001
$ws.Cells.Item(2, 9).Value = 676.92
$ws.Cells.Item(2, 11).Value = 676.92
$ws.Cells.Item(2, 13).Value = -563.92

$ws.Cells.Item(16, 8).Value = 17552.75
$ws.Cells.Item(16, 9).Value = 10070.333
$ws.Cells.Item(16, 10).Value = 40000
$ws.Cells.Item(16, 11).Value = 10070.333
$ws.Cells.Item(16, 12).Value = 40000
$ws.Cells.Item(16, 13).Value = -9783.333000000001
$ws.Cells.Item(16, 14).Value = -40574

$ws.Cells.Item(32, 8).Value = 15402.933
$ws.Cells.Item(32, 9).Value = 14366.204
$ws.Cells.Item(32, 11).Value = 14366.204
$ws.Cells.Item(32, 13).Value = -14079.204

$ws.Cells.Item(116, 8).Value = 693.2593000000001
$ws.Cells.Item(116, 9).Value = 676.92
$ws.Cells.Item(116, 11).Value = 676.92
$ws.Cells.Item(116, 13).Value = 1617.08

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 693.2593000000001
$ws.Cells.Item(3, 9).Value = 676.92
$ws.Cells.Item(3, 11).Value = 676.92
$ws.Cells.Item(3, 13).Value = -562.92

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 63493690
$ws.Cells.Item(132, 9).Value = 63493690
$ws.Cells.Item(132, 11).Value = 190481070
$ws.Cells.Item(132, 13).Value = -190478540

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(82, 8).Value = 5249.875
$ws.Cells.Item(82, 10).Value = 6400
$ws.Cells.Item(82, 12).Value = 19200
$ws.Cells.Item(82, 14).Value = -20012

$ws.Cells.Item(85, 8).Value = 5249.875
$ws.Cells.Item(85, 10).Value = 6400
$ws.Cells.Item(85, 12).Value = 19200
$ws.Cells.Item(85, 14).Value = -22008

$ws.Cells.Item(113, 8).Value = 684.0714
$ws.Cells.Item(113, 9).Value = 608
$ws.Cells.Item(113, 10).Value = 704.8182
$ws.Cells.Item(113, 11).Value = 1824
$ws.Cells.Item(113, 12).Value = 2114.4546
$ws.Cells.Item(113, 13).Value = 346
$ws.Cells.Item(113, 14).Value = -6454.4546

$ws.Cells.Item(119, 8).Value = 15096.6
$ws.Cells.Item(119, 9).Value = 12163.667
$ws.Cells.Item(119, 10).Value = 19496
$ws.Cells.Item(119, 11).Value = 36491.001
$ws.Cells.Item(119, 12).Value = 58488
$ws.Cells.Item(119, 13).Value = -31653.001
$ws.Cells.Item(119, 14).Value = -68164

$ws.Cells.Item(121, 8).Value = 1119.909
$ws.Cells.Item(121, 10).Value = 1443.375
$ws.Cells.Item(121, 12).Value = 4330.125
$ws.Cells.Item(121, 14).Value = -6950.125

$ws.Cells.Item(124, 8).Value = 3683
$ws.Cells.Item(124, 9).Value = 4666
$ws.Cells.Item(124, 10).Value = 2700
$ws.Cells.Item(124, 11).Value = 13998
$ws.Cells.Item(124, 12).Value = 8100
$ws.Cells.Item(124, 13).Value = -9088
$ws.Cells.Item(124, 14).Value = -17920

$ws.Cells.Item(126, 8).Value = 1000
$ws.Cells.Item(126, 9).Value = 1000
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 3000
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).ClearContents()
$ws.Cells.Item(126, 14).Value = 1940

$ws.Cells.Item(131, 8).Value = 2955.76
$ws.Cells.Item(131, 10).Value = 3023.1765
$ws.Cells.Item(131, 12).Value = 9069.529500000001
$ws.Cells.Item(131, 14).Value = -19149.5295

$ws.Cells.Item(137, 8).Value = 2438.6365
$ws.Cells.Item(137, 9).Value = 2066
$ws.Cells.Item(137, 10).Value = 2578.375
$ws.Cells.Item(137, 11).Value = 6198
$ws.Cells.Item(137, 12).Value = 7735.125
$ws.Cells.Item(137, 13).Value = -1098
$ws.Cells.Item(137, 14).Value = -17935.125

$ws.Cells.Item(140, 8).Value = 4353.8
$ws.Cells.Item(140, 9).Value = 4170.8887
$ws.Cells.Item(140, 11).Value = 12512.6661
$ws.Cells.Item(140, 13).Value = -7332.666100000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 361.42856
$ws.Cells.Item(2, 9).Value = 132.2
$ws.Cells.Item(2, 10).Value = 488.77777
$ws.Cells.Item(2, 11).Value = 132.2
$ws.Cells.Item(2, 12).Value = 488.77777
$ws.Cells.Item(2, 13).Value = -19.19999999999999
$ws.Cells.Item(2, 14).Value = -714.7777699999999

$ws.Cells.Item(21, 8).Value = 21000
$ws.Cells.Item(21, 9).Value = 21000
$ws.Cells.Item(21, 11).Value = 21000
$ws.Cells.Item(21, 13).Value = -20827

$ws.Cells.Item(30, 8).Value = 21000
$ws.Cells.Item(30, 9).Value = 21000
$ws.Cells.Item(30, 11).Value = 21000
$ws.Cells.Item(30, 13).Value = -20895

$ws.Cells.Item(96, 8).Value = 39499.5
$ws.Cells.Item(96, 10).Value = 39499.5
$ws.Cells.Item(96, 12).Value = 39499.5
$ws.Cells.Item(96, 14).Value = -44991.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 937.8570999999999
$ws.Cells.Item(22, 9).Value = 773.1
$ws.Cells.Item(22, 10).Value = 1349.75
$ws.Cells.Item(22, 11).Value = 773.1
$ws.Cells.Item(22, 12).Value = 1349.75
$ws.Cells.Item(22, 13).Value = -478.1
$ws.Cells.Item(22, 14).Value = -1939.75

$ws.Cells.Item(27, 8).Value = 937.8570999999999
$ws.Cells.Item(27, 9).Value = 773.1
$ws.Cells.Item(27, 10).Value = 1349.75
$ws.Cells.Item(27, 11).Value = 773.1
$ws.Cells.Item(27, 12).Value = 1349.75
$ws.Cells.Item(27, 13).Value = -666.1
$ws.Cells.Item(27, 14).Value = -1563.75

$ws.Cells.Item(122, 8).Value = 11926.214
$ws.Cells.Item(122, 9).Value = 5001
$ws.Cells.Item(122, 10).Value = 13814.909
$ws.Cells.Item(122, 11).Value = 15003
$ws.Cells.Item(122, 12).Value = 41444.727
$ws.Cells.Item(122, 13).Value = -12553
$ws.Cells.Item(122, 14).Value = -46344.727

$ws.Cells.Item(132, 8).Value = 2798.9
$ws.Cells.Item(132, 9).Value = 2794.027
$ws.Cells.Item(132, 11).Value = 8382.081
$ws.Cells.Item(132, 13).Value = -5852.081

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 2496.6667
$ws.Cells.Item(62, 9).Value = 2245
$ws.Cells.Item(62, 11).Value = 2245
$ws.Cells.Item(62, 13).Value = -1621

$ws.Cells.Item(65, 8).Value = 2496.6667
$ws.Cells.Item(65, 9).Value = 2245
$ws.Cells.Item(65, 11).Value = 11225
$ws.Cells.Item(65, 13).Value = -8105

$ws.Cells.Item(107, 8).Value = 858.86664
$ws.Cells.Item(107, 10).Value = 970
$ws.Cells.Item(107, 12).Value = 2910
$ws.Cells.Item(107, 14).Value = -6750
